# Update the FHIR StructureDefinition workbook from the Alvearie/IBM
# publication values to the LinuxForHealth publication values.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet ----
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-category"
# Version
$meta.Range("B3").Value = "8.0.0"
# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---- Elements sheet ----
$elements = $wb.Worksheets.Item("Elements")

# The root Extension row no longer carries the ele-1/ext-1 constraint text
# in the "Constraint(s)" column (AI2).
$elements.Range("AI2").Value = ""

# The Extension.url row's "Fixed Value" column (Q5) mirrors the new URL.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-category"
